# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Wed Oct 23 04:51:03 UTC 2024 with GitHub Actions".
# For each affected row (2-51) update the Price (column D, when present)
# and Volume(1h) (column E) cells with the newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    # Plain assignment is fine for values Excel will not mistake for a
    # number (ids like "67.188.62" or the subscript-formatted price in
    # row 30), but pure-number-looking strings ("1.00", "7.79", ...) get
    # silently coerced to a Double, which drops an exact trailing zero.
    # Forcing a Text number format for the write keeps it a literal
    # string; restoring the "Normal" style right after avoids leaving a
    # stray number format behind on the cell. None of the prices touched
    # here are actually "0", so a failed/non-numeric parse (which this
    # host surfaces as 0) is an unambiguous signal to leave it alone.
    $trimmed = $text.Trim()
    $looksNumeric = $false
    if ($trimmed.Length -gt 0) {
        $n = [double]$trimmed
        if ($n -ne 0) { $looksNumeric = $true }
    }
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Cells.Item(2, 4) "67.188.62"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.32%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.616.27"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.85%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.00"
Set-TextValue $ws.Cells.Item(4, 5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(5, 4) "590.38"
Set-TextValue $ws.Cells.Item(5, 5) "  -1.27%  "
Set-TextValue $ws.Cells.Item(6, 4) "165.76"
Set-TextValue $ws.Cells.Item(6, 5) "  -0.69%  "
Set-TextValue $ws.Cells.Item(7, 5) "  +0.01%  "
Set-TextValue $ws.Cells.Item(8, 5) "  -2.14%  "
Set-TextValue $ws.Cells.Item(9, 4) "2.615.20"
Set-TextValue $ws.Cells.Item(9, 5) "  -0.87%  "
Set-TextValue $ws.Cells.Item(10, 5) "  -4.04%  "
Set-TextValue $ws.Cells.Item(11, 5) "  +1.52%  "
Set-TextValue $ws.Cells.Item(12, 5) "  -0.27%  "
Set-TextValue $ws.Cells.Item(13, 5) "  -0.45%  "
Set-TextValue $ws.Cells.Item(14, 4) "27.34"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.12%  "
Set-TextValue $ws.Cells.Item(15, 4) "3.091.63"
Set-TextValue $ws.Cells.Item(15, 5) "  -1.00%  "
Set-TextValue $ws.Cells.Item(16, 5) "  -2.26%  "
Set-TextValue $ws.Cells.Item(17, 4) "67.237.81"
Set-TextValue $ws.Cells.Item(17, 5) "  -0.16%  "
Set-TextValue $ws.Cells.Item(18, 4) "2.640.42"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.51%  "
Set-TextValue $ws.Cells.Item(19, 4) "11.77"
Set-TextValue $ws.Cells.Item(19, 5) "  -0.74%  "
Set-TextValue $ws.Cells.Item(20, 4) "7.79"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.57%  "
Set-TextValue $ws.Cells.Item(21, 4) "354.80"
Set-TextValue $ws.Cells.Item(21, 5) "  -2.18%  "
Set-TextValue $ws.Cells.Item(22, 5) "  -2.89%  "
Set-TextValue $ws.Cells.Item(23, 5) "  -2.87%  "
Set-TextValue $ws.Cells.Item(24, 4) "10.53"
Set-TextValue $ws.Cells.Item(24, 5) "  -3.71%  "
Set-TextValue $ws.Cells.Item(25, 5) "  -0.06%  "
Set-TextValue $ws.Cells.Item(26, 5) "  -4.34%  "
Set-TextValue $ws.Cells.Item(27, 5) "  -2.36%  "
Set-TextValue $ws.Cells.Item(28, 5) "  -1.20%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.00"
Set-TextValue $ws.Cells.Item(29, 5) "  -0.01%  "
Set-TextValue $ws.Cells.Item(30, 4) "0.0₃0997"
Set-TextValue $ws.Cells.Item(30, 5) "  -2.49%  "
Set-TextValue $ws.Cells.Item(31, 4) "543.19"
Set-TextValue $ws.Cells.Item(31, 5) "  -1.84%  "
Set-TextValue $ws.Cells.Item(32, 4) "7.86"
Set-TextValue $ws.Cells.Item(32, 5) "  -2.27%  "
Set-TextValue $ws.Cells.Item(33, 5) "  -3.68%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.87"
Set-TextValue $ws.Cells.Item(34, 5) "  -2.63%  "
Set-TextValue $ws.Cells.Item(35, 5) "  +0.31%  "
Set-TextValue $ws.Cells.Item(36, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.48"
Set-TextValue $ws.Cells.Item(37, 5) "  -3.46%  "
Set-TextValue $ws.Cells.Item(38, 4) "159.15"
Set-TextValue $ws.Cells.Item(38, 5) "  +0.91%  "
Set-TextValue $ws.Cells.Item(39, 4) "18.91"
Set-TextValue $ws.Cells.Item(39, 5) "  -2.44%  "
Set-TextValue $ws.Cells.Item(40, 5) "  -2.14%  "
Set-TextValue $ws.Cells.Item(41, 4) "18.24"
Set-TextValue $ws.Cells.Item(41, 5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(42, 5) "  -1.31%  "
Set-TextValue $ws.Cells.Item(43, 4) "5.14"
Set-TextValue $ws.Cells.Item(43, 5) "  -2.20%  "
Set-TextValue $ws.Cells.Item(44, 5) "  +0.04%  "
Set-TextValue $ws.Cells.Item(45, 5) "  -4.37%  "
Set-TextValue $ws.Cells.Item(46, 5) "  -0.95%  "
Set-TextValue $ws.Cells.Item(47, 4) "151.58"
Set-TextValue $ws.Cells.Item(47, 5) "  -1.19%  "
Set-TextValue $ws.Cells.Item(48, 5) "  -3.36%  "
Set-TextValue $ws.Cells.Item(49, 4) "3.76"
Set-TextValue $ws.Cells.Item(49, 5) "  -2.95%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.70"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.17%  "
Set-TextValue $ws.Cells.Item(51, 5) "  -1.29%  "
